$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title ("Play Acorn Pixie Slot Free: Enchanted World of Fairies").
$metaRange = $d.Content
$metaRange.Find.Execute("Meta description") | Out-Null
$metaRange.Paragraphs(1).Range.Delete()

# 2. Insert a new bold paragraph with the title text right before the final
#    paragraph (the one that currently holds the italic "Prompt for DALLE..." text).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($count)
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Acorn Pixie Slot Free: Enchanted World of Fairies</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($titleXml)

# 3. Swap the DALLE image prompt text in the final (italic) paragraph for the
#    meta-description copy, keeping the existing italic run formatting intact.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalPara.Range.Find.Execute("Prompt for DALLE: Create a feature image for Acorn Pixie that showcases a happy Maya warrior in cartoon style with glasses. The image should incorporate elements from the game such as fairies, elves, flowers, and mushrooms, with a mystical forest backdrop. The Maya warrior can be holding an acorn or surrounded by them, with the game's logo appearing somewhere in the image. Make the image bright, colorful, and enticing to capture the attention of potential players.", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "Experience the mystical feel of Acorn Pixie with visually stunning designs and dynamic gameplay. Play free and unlock Cluster Wild and Acorn Pixie Bonus features.", 2)
